{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the Review_228 diff:\n//  - para 0: date 19.06.24 -> 18.06.24\n//  - para 1: title replaced\n//  - para 2: body replaced\n//  - para 3: body replaced\n//  - para 4: body replaced with new arXiv link (used to be \"Fastfood\" paragraph)\n//  - para 5 & 6 removed (old \"\u05d1\u05e7\u05d9\u05e6\u05d5\u05e8\" paragraph and old arXiv link paragraph)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst newTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 18.06.24:\u26a1\ufe0f\ud83d\ude80\",\n  \"Helping or Herding? Reward Model Ensembles Mitigate but do not Eliminate REWARD HACKING\",\n  \" \u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d4\u05d6\u05d5 \u05de\u05de\u05e9\u05d9\u05db\u05d4 \u05d0\u05ea \u05e7\u05d5 \u05d4\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea \u05e9\u05db\u05ea\u05d1\u05ea\u05d9 \u05d1\u05e0\u05d5\u05e9\u05d0 RLHF. \u05db\u05de\u05d5 \u05e9\u05d0\u05ea\u05dd \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 \u05d1-RLHF \u05de\u05db\u05d9\u05dc\u05d4 \u05e9\u05e0\u05d9 \u05d0\u05d9\u05d1\u05e8\u05d9\u05dd: \u05d4\u05d0\u05d9\u05d1\u05e8 \u05e9\u05de\u05e0\u05e1\u05d4 \u05dc\u05de\u05e7\u05e1\u05dd \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc (reward) \u05d5\u05d4\u05d0\u05d9\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05de\u05e0\u05e1\u05d4 \u05dc\u05e9\u05de\u05d5\u05e8 \u05d0\u05ea \u05de\u05d5\u05d3\u05dc \u05d4\u05e9\u05e4\u05d4 \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d9\u05d5\u05d1 (\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05e1\u05d5\u05e4\u05d9) \u05e7\u05e8\u05d5\u05d1 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4-RLHF \u05de\u05de\u05e0\u05d5. \u05d1\u05e2\u05d1\u05e8 \u05d9\u05e6\u05d0\u05d5 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05d4\u05e6\u05d9\u05e2\u05d5 \u05dc\u05d0\u05de\u05df \u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d5\u05d0\u05d6 \u05dc\u05de\u05e6\u05e2 (\u05d0\u05d5 \u05dc\u05e7\u05d7\u05ea \u05de\u05e7\u05e1\u05d9\u05de\u05d5\u05dd) \u05e9\u05dc \u05db\u05dc \u05d4-rewards \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05e9\u05d0\u05dc\u05d4 \u05d5\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e0\u05ea\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \u05d6\u05d4 \u05dc\u05d8\u05e2\u05e0\u05ea\u05dd \u05de\u05e7\u05d8\u05d9\u05df \u05d0\u05ea \u05d4\u05e1\u05d9\u05db\u05d5\u05d9 \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d1-RLHF \u05d9\u05d1\u05e6\u05e2 reward hacking \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05ea\u05db\u05e0\u05e1 \u05dc\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05de\u05de\u05e7\u05e1\u05dd \u05ea\u05d2\u05de\u05d5\u05dc \u05d0\u05da \u05d1\u05e4\u05e2\u05d5\u05dc \u05de\u05d2\u05e0\u05e8\u05d8 \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05e8\u05d5\u05e2\u05d4.\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05d8\u05d5\u05e2\u05df \u05e9\u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5 \u05d0\u05d9\u05e0\u05d4 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05ea \u05db\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 \u05e9\u05d0\u05d9\u05ea\u05d4 \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 reward (\u05db\u05dc\u05d5\u05de\u05e8 Bradley-Terry) \u05d2\u05d5\u05e8\u05de\u05ea \u05dc\u05db\u05da \u05e9\u05db\u05dc \u05e9\u05e0\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05e9\u05d5\u05e0\u05d9\u05dd \u05e8\u05e7 \u05d1\u05e7\u05d1\u05d5\u05e2 \u05e9\u05ea\u05dc\u05d5\u05d9 \u05e8\u05e7 \u05d1\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x \u05d9\u05e7\u05d1\u05dc\u05d5 \u05d0\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d4\u05e2\u05e8\u05da \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1. \u05d1\u05e4\u05d5\u05e2\u05dc \u05d6\u05d4 \u05d0\u05d5\u05de\u05e8 \u05db\u05d9 \u05dc\u05db\u05dc \u05dc\u05e2\u05e8\u05db\u05d9 \u05d4- reward\u05bf, \u05d4\u05de\u05d5\u05e4\u05e7\u05d9\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9, \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05de\u05d5\u05e6\u05e2\u05d9\u05dd \u05d5\u05d1\u05e4\u05d5\u05e2\u05dc \u05d4\u05d1\u05d7\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d4\u05de\u05e7\u05e1\u05d9\u05de\u05dc\u05d9 \u05d0\u05d5 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05de\u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 \u05db\u05d0\u05dc\u05d5 \u05e2\u05e9\u05d5\u05d9\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05dc\u05d0 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05ea (\u05db\u05de\u05d5 \u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05ea\u05e4\u05d5\u05d6\u05d9\u05dd \u05d5\u05e2\u05d2\u05d1\u05e0\u05d9\u05d4). \u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05e2\u05dd \u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05d1\u05d0\u05d4 \u05f4\u05dc\u05e8\u05e1\u05df\u05f4 \u05d0\u05ea \u05d4\u05e7\u05d1\u05d5\u05e2 \u05d6\u05d4 \u05e9\u05ea\u05dc\u05d5\u05d9 \u05e8\u05e7 \u05d1\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05d1\u05db\u05da \u05f4\u05dc\u05e1\u05db\u05e0\u05e8\u05df\u05f4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05e9\u05d5\u05e0\u05d9\u05dd.\",\n  \"https://arxiv.org/abs/2312.09244\",\n];\n\nconst items = paragraphs.items;\n\n// Replace the text of the first five paragraphs (indices 0-4) in place.\nfor (let i = 0; i < newTexts.length; i++) {\n  items[i].insertText(newTexts[i], \"Replace\");\n}\n\n// Remove the two trailing paragraphs that no longer exist in the new version\n// (the old \"Fastfood\"-era wrap-up sentence and the old arXiv link paragraph).\nfor (let i = items.length - 1; i >= newTexts.length; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the Review_228 diff:\n#  - para 1: date 19.06.24 -> 18.06.24\n#  - para 2: title replaced\n#  - para 3: body replaced\n#  - para 4: body replaced\n#  - para 5: body replaced with new arXiv link (used to be \"Fastfood\" paragraph)\n#  - para 6 & 7 removed (old \"\u05d1\u05e7\u05d9\u05e6\u05d5\u05e8\" paragraph and old arXiv link paragraph)\n\n$d = $word.ActiveDocument\n\n$d.Paragraphs.Item(1).Range.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 18.06.24:\u26a1\ufe0f\ud83d\ude80\"\n$d.Paragraphs.Item(2).Range.Text = \"Helping or Herding? Reward Model Ensembles Mitigate but do not Eliminate REWARD HACKING\"\n$d.Paragraphs.Item(3).Range.Text = \" \u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05d4\u05d6\u05d5 \u05de\u05de\u05e9\u05d9\u05db\u05d4 \u05d0\u05ea \u05e7\u05d5 \u05d4\u05e1\u05e7\u05d9\u05e8\u05d5\u05ea \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea \u05e9\u05db\u05ea\u05d1\u05ea\u05d9 \u05d1\u05e0\u05d5\u05e9\u05d0 RLHF. \u05db\u05de\u05d5 \u05e9\u05d0\u05ea\u05dd \u05d6\u05d5\u05db\u05e8\u05d9\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 \u05d1-RLHF \u05de\u05db\u05d9\u05dc\u05d4 \u05e9\u05e0\u05d9 \u05d0\u05d9\u05d1\u05e8\u05d9\u05dd: \u05d4\u05d0\u05d9\u05d1\u05e8 \u05e9\u05de\u05e0\u05e1\u05d4 \u05dc\u05de\u05e7\u05e1\u05dd \u05d0\u05ea \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc (reward) \u05d5\u05d4\u05d0\u05d9\u05d1\u05e8 \u05d4\u05e9\u05e0\u05d9 \u05de\u05e0\u05e1\u05d4 \u05dc\u05e9\u05de\u05d5\u05e8 \u05d0\u05ea \u05de\u05d5\u05d3\u05dc \u05d4\u05e9\u05e4\u05d4 \u05d0\u05d7\u05e8\u05d9 \u05d8\u05d9\u05d5\u05d1 (\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05e1\u05d5\u05e4\u05d9) \u05e7\u05e8\u05d5\u05d1 \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05de\u05ea\u05d7\u05d9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4-RLHF \u05de\u05de\u05e0\u05d5. \u05d1\u05e2\u05d1\u05e8 \u05d9\u05e6\u05d0\u05d5 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05d4\u05e6\u05d9\u05e2\u05d5 \u05dc\u05d0\u05de\u05df \u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05d5\u05d0\u05d6 \u05dc\u05de\u05e6\u05e2 (\u05d0\u05d5 \u05dc\u05e7\u05d7\u05ea \u05de\u05e7\u05e1\u05d9\u05de\u05d5\u05dd) \u05e9\u05dc \u05db\u05dc \u05d4-rewards \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05e9\u05d0\u05dc\u05d4 \u05d5\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e0\u05ea\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4. \u05d6\u05d4 \u05dc\u05d8\u05e2\u05e0\u05ea\u05dd \u05de\u05e7\u05d8\u05d9\u05df \u05d0\u05ea \u05d4\u05e1\u05d9\u05db\u05d5\u05d9 \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05d1-RLHF \u05d9\u05d1\u05e6\u05e2 reward hacking \u05db\u05dc\u05d5\u05de\u05e8 \u05d9\u05ea\u05db\u05e0\u05e1 \u05dc\u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d4\u05de\u05de\u05e7\u05e1\u05dd \u05ea\u05d2\u05de\u05d5\u05dc \u05d0\u05da \u05d1\u05e4\u05e2\u05d5\u05dc \u05de\u05d2\u05e0\u05e8\u05d8 \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05d1\u05d0\u05d9\u05db\u05d5\u05ea \u05d2\u05e8\u05d5\u05e2\u05d4.\"\n$d.Paragraphs.Item(4).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05e0\u05e1\u05e7\u05d5\u05e8 \u05d4\u05d9\u05d5\u05dd \u05d8\u05d5\u05e2\u05df \u05e9\u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5 \u05d0\u05d9\u05e0\u05d4 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05ea \u05db\u05d9 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1 \u05e9\u05d0\u05d9\u05ea\u05d4 \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 reward (\u05db\u05dc\u05d5\u05de\u05e8 Bradley-Terry) \u05d2\u05d5\u05e8\u05de\u05ea \u05dc\u05db\u05da \u05e9\u05db\u05dc \u05e9\u05e0\u05d9 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05e9\u05d5\u05e0\u05d9\u05dd \u05e8\u05e7 \u05d1\u05e7\u05d1\u05d5\u05e2 \u05e9\u05ea\u05dc\u05d5\u05d9 \u05e8\u05e7 \u05d1\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x \u05d9\u05e7\u05d1\u05dc\u05d5 \u05d0\u05ea \u05d0\u05d5\u05ea\u05d5 \u05d4\u05e2\u05e8\u05da \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05dc\u05d5\u05e1. \u05d1\u05e4\u05d5\u05e2\u05dc \u05d6\u05d4 \u05d0\u05d5\u05de\u05e8 \u05db\u05d9 \u05dc\u05db\u05dc \u05dc\u05e2\u05e8\u05db\u05d9 \u05d4- reward\u05bf, \u05d4\u05de\u05d5\u05e4\u05e7\u05d9\u05dd \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05de\u05d5\u05d3\u05dc\u05d9, \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05de\u05d5\u05e6\u05e2\u05d9\u05dd \u05d5\u05d1\u05e4\u05d5\u05e2\u05dc \u05d4\u05d1\u05d7\u05d9\u05e8\u05d4 \u05e9\u05dc \u05d4\u05de\u05e7\u05e1\u05d9\u05de\u05dc\u05d9 \u05d0\u05d5 \u05d4\u05de\u05de\u05d5\u05e6\u05e2 \u05de\u05db\u05de\u05d4 \u05de\u05d5\u05d3\u05dc\u05d9 \u05db\u05d0\u05dc\u05d5 \u05e2\u05e9\u05d5\u05d9\u05d4 \u05dc\u05d4\u05d9\u05d5\u05ea \u05dc\u05d0 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05dc\u05d9\u05ea (\u05db\u05de\u05d5 \u05de\u05de\u05d5\u05e6\u05e2 \u05e9\u05dc \u05ea\u05e4\u05d5\u05d6\u05d9\u05dd \u05d5\u05e2\u05d2\u05d1\u05e0\u05d9\u05d4). \u05d0\u05d6 \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d0\u05de\u05df \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05e2\u05dd \u05e8\u05d2\u05d5\u05dc\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e9\u05d1\u05d0\u05d4 \u05f4\u05dc\u05e8\u05e1\u05df\u05f4 \u05d0\u05ea \u05d4\u05e7\u05d1\u05d5\u05e2 \u05d6\u05d4 \u05e9\u05ea\u05dc\u05d5\u05d9 \u05e8\u05e7 \u05d1\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 \u05d5\u05d1\u05db\u05da \u05f4\u05dc\u05e1\u05db\u05e0\u05e8\u05df\u05f4 \u05de\u05d5\u05d3\u05dc\u05d9 reward \u05e9\u05d5\u05e0\u05d9\u05dd.\"\n$d.Paragraphs.Item(5).Range.Text = \"https://arxiv.org/abs/2312.09244\"\n\n# Remove the two now-obsolete trailing paragraphs (delete from the end so the\n# indices of the paragraphs we still need stay stable).\n$d.Paragraphs.Item(7).Range.Delete()\n$d.Paragraphs.Item(6).Range.Delete()\n"}
